$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'66.702.00"
$ws.Range("E2").Value = "  -0.74%  "

# Row 3
$ws.Range("D3").Value = "'3.792.41"
$ws.Range("E3").Value = "  -1.31%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'431.68"
$ws.Range("E5").Value = "  +4.34%  "

# Row 6
$ws.Range("D6").Value = "'140.77"
$ws.Range("E6").Value = "  +6.10%  "

# Row 7
$ws.Range("D7").Value = "'0.627"
$ws.Range("E7").Value = "  +1.86%  "

# Row 8
$ws.Range("E8").Value = "  -0.13%  "

# Row 9
$ws.Range("E9").Value = "  -0.95%  "

# Row 10
$ws.Range("D10").Value = "'0.154"
$ws.Range("E10").Value = "  -10.19%  "

# Row 11
$ws.Range("E11").Value = "  -16.90%  "

# Row 12
$ws.Range("D12").Value = "'42.80"
$ws.Range("E12").Value = "  +3.96%  "

# Row 13
$ws.Range("D13").Value = "'10.47"
$ws.Range("E13").Value = "  +3.88%  "

# Row 14
$ws.Range("D14").Value = "'4.386.98"
$ws.Range("E14").Value = "  -1.03%  "

# Row 15
$ws.Range("D15").Value = "'15.06"
$ws.Range("E15").Value = "  +1.28%  "

# Row 16
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "'0.138"
$ws.Range("E16").Value = "  -0.19%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'3.793.52"
$ws.Range("E17").Value = "  -2.14%  "

# Row 18
$ws.Range("D18").Value = "'19.97"
$ws.Range("E18").Value = "  +2.06%  "

# Row 19
$ws.Range("D19").Value = "'1.13"
$ws.Range("E19").Value = "  +5.50%  "

# Row 20
$ws.Range("D20").Value = "'66.749.94"
$ws.Range("E20").Value = "  -1.10%  "

# Row 21
$ws.Range("D21").Value = "'410.23"
$ws.Range("E21").Value = "  -1.56%  "

# Row 22
$ws.Range("D22").Value = "'14.74"
$ws.Range("E22").Value = "  -0.24%  "

# Row 23
$ws.Range("D23").Value = "'3.28"
$ws.Range("E23").Value = "  +6.33%  "

# Row 24
$ws.Range("D24").Value = "'85.30"
$ws.Range("E24").Value = "  -0.90%  "

# Row 25
$ws.Range("D25").Value = "'36.79"
$ws.Range("E25").Value = "  -0.43%  "

# Row 26
$ws.Range("D26").Value = "'3.34"
$ws.Range("E26").Value = "  +5.98%  "

# Row 27
$ws.Range("E27").Value = "  -2.42%  "

# Row 28
$ws.Range("D28").Value = "'9.68"
$ws.Range("E28").Value = "  +33.31%  "

# Row 29
$ws.Range("D29").Value = "'9.78"
$ws.Range("E29").Value = "  +2.38%  "

# Row 30
$ws.Range("D30").Value = "'726.31"
$ws.Range("E30").Value = "  +5.48%  "

# Row 31
$ws.Range("D31").Value = "'13.86"
$ws.Range("E31").Value = "  +10.34%  "

# Row 32
$ws.Range("D32").Value = "'0.135"
$ws.Range("E32").Value = "  +9.77%  "

# Row 33
$ws.Range("D33").Value = "'2.76"
$ws.Range("E33").Value = "  +0.66%  "

# Row 34
$ws.Range("D34").Value = "'41.65"
$ws.Range("E34").Value = "  +6.21%  "

# Row 35
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.09%  "

# Row 36
$ws.Range("E36").Value = "  +28.85%  "

# Row 37
$ws.Range("E37").Value = "  -1.89%  "

# Row 38
$ws.Range("E38").Value = "  +0.56%  "

# Row 39
$ws.Range("E39").Value = "  +2.64%  "

# Row 40
$ws.Range("D40").Value = "'2.73"
$ws.Range("E40").Value = "  +39.74%  "

# Row 41
$ws.Range("D41").Value = "'2.97"
$ws.Range("E41").Value = "  -3.55%  "

# Row 42
$ws.Range("D42").Value = "0.0₃0686"
$ws.Range("E42").Value = "  -13.87%  "

# Row 43
$ws.Range("E43").Value = "  +3.13%  "

# Row 44
$ws.Range("E44").Value = "  +0.57%  "

# Row 45
$ws.Range("D45").Value = "'3.22"
$ws.Range("E45").Value = "  +0.61%  "

# Row 46
$ws.Range("E46").Value = "  +8.46%  "

# Row 47
$ws.Range("B47").Value = "LidoDAOToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D47").Value = "'3.33"
$ws.Range("E47").Value = "  -0.35%  "

# Row 48
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "'2.69"
$ws.Range("E48").Value = "  +3.67%  "

# Row 49
$ws.Range("E49").Value = "  -0.94%  "

# Row 50
$ws.Range("D50").Value = "'142.29"
$ws.Range("E50").Value = "  -4.11%  "

# Row 51
$ws.Range("E51").Value = "  -1.18%  "
